$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Calc reference style switched to R1C1 (as recorded in the workbook's calcPr)
$excel.ReferenceStyle = "xlR1C1"

# Continue the daily chinups/pushups/eggs-style tracker with another week
# of entries (rows 9-15), matching the date-formatted column A used by
# the existing rows.
$newRows = @(
    @(44223, 0, 0),
    @(44224, 0, 20),
    @(44225, 0, 0),
    @(44226, 0, 0),
    @(44227, 0, 0),
    @(44228, 0, 0),
    @(44229, 0, 0)
)

$startRow = 9
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]

    # Carry the date number format down from the row above (column A only)
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Leave the selection where the user's cursor ended up after typing the
# last entry.
$ws.Range("E15").Select() | Out-Null
